# Scheduled-runner refresh of live market/profit figures across the
# "Halicarnassus Profits" leve-crafting workbook. Each block below
# rewrites the computed price/profit columns (H:N) for one leve row,
# mirroring the upstream data refresh described in the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 304.75
$ws.Range("I2").Value = 110
$ws.Range("K2").Value = 110
$ws.Range("M2").Value = 3

$ws.Range("H11").Value = 97.90909000000001
$ws.Range("I11").Value = 97.90909000000001
$ws.Range("K11").Value = 97.90909000000001
$ws.Range("M11").Value = 42.09090999999999

$ws.Range("H55").Value = 756.75
$ws.Range("I55").Value = 788.625
$ws.Range("J55").Value = 735.5
$ws.Range("K55").Value = 788.625
$ws.Range("L55").Value = 735.5
$ws.Range("M55").Value = -574.625
$ws.Range("N55").Value = -1163.5

$ws.Range("H129").Value = 1299
$ws.Range("I129").Value = 1062.3334
$ws.Range("K129").Value = 3187.0002
$ws.Range("M129").Value = 1812.9998

$ws.Range("H137").Value = 2210.3333
$ws.Range("I137").Value = 544.8
$ws.Range("J137").Value = 3400
$ws.Range("K137").Value = 1634.4
$ws.Range("L137").Value = 10200
$ws.Range("M137").Value = 915.6000000000001
$ws.Range("N137").Value = -15300

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 7225.875
$ws.Range("I102").Value = 4265.3335
$ws.Range("J102").Value = 9002.200000000001
$ws.Range("K102").Value = 4265.3335
$ws.Range("L102").Value = 9002.200000000001
$ws.Range("M102").Value = -2643.3335
$ws.Range("N102").Value = -12246.2

$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("M122").Value = -3550

$ws.Range("H132").Value = 4388.385
$ws.Range("I132").Value = 5169.4
$ws.Range("K132").Value = 15508.2
$ws.Range("M132").Value = -12978.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").ClearContents()
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = 0

$ws.Range("H26").Value = 21388.834
$ws.Range("I26").Value = 21388.834
$ws.Range("K26").Value = 21388.834
$ws.Range("M26").Value = -21096.834

$ws.Range("H86").Value = 6033.222
$ws.Range("I86").Value = 1766.6666
$ws.Range("K86").Value = 1766.6666
$ws.Range("M86").Value = -643.6666

$ws.Range("H89").Value = 6033.222
$ws.Range("I89").Value = 1766.6666
$ws.Range("K89").Value = 8833.333000000001
$ws.Range("M89").Value = -3217.333000000001

$ws.Range("H96").Value = 16900.75
$ws.Range("I96").Value = 16900.75
$ws.Range("K96").Value = 16900.75
$ws.Range("M96").Value = -14154.75

$ws.Range("H99").Value = 2009.2222
$ws.Range("I99").Value = 1440.5714
$ws.Range("K99").Value = 1440.5714
$ws.Range("M99").Value = 57.42859999999996

$ws.Range("H107").Value = 7053.143
$ws.Range("I107").Value = 7220.5
$ws.Range("K107").Value = 7220.5
$ws.Range("M107").Value = -5300.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 8877.75
$ws.Range("J62").Value = 10005.5
$ws.Range("L62").Value = 10005.5
$ws.Range("N62").Value = -11253.5

$ws.Range("H65").Value = 8877.75
$ws.Range("J65").Value = 10005.5
$ws.Range("L65").Value = 50027.5
$ws.Range("N65").Value = -56267.5

$ws.Range("H68").Value = 46951.2
$ws.Range("J68").Value = 46951.2
$ws.Range("L68").Value = 46951.2
$ws.Range("N68").Value = -48449.2

$ws.Range("H69").Value = 13536.4
$ws.Range("I69").Value = 10670.5
$ws.Range("J69").Value = 25000
$ws.Range("K69").Value = 10670.5
$ws.Range("L69").Value = 25000
$ws.Range("M69").Value = -9921.5
$ws.Range("N69").Value = -26498

$ws.Range("H71").Value = 46951.2
$ws.Range("J71").Value = 46951.2
$ws.Range("L71").Value = 140853.6
$ws.Range("N71").Value = -148341.6

$ws.Range("H72").Value = 13536.4
$ws.Range("I72").Value = 10670.5
$ws.Range("J72").Value = 25000
$ws.Range("K72").Value = 32011.5
$ws.Range("L72").Value = 75000
$ws.Range("M72").Value = -28267.5
$ws.Range("N72").Value = -82488

$ws.Range("H80").Value = 20000
$ws.Range("J80").Value = 20000
$ws.Range("L80").Value = 20000
$ws.Range("N80").Value = -22246

$ws.Range("H83").Value = 20000
$ws.Range("J83").Value = 20000
$ws.Range("L83").Value = 60000
$ws.Range("N83").Value = -71232

$ws.Range("H93").Value = 4951.6665
$ws.Range("I93").Value = 5407
$ws.Range("J93").Value = 4724
$ws.Range("K93").Value = 5407
$ws.Range("L93").Value = 4724
$ws.Range("M93").Value = -3535
$ws.Range("N93").Value = -8468

$ws.Range("H103").Value = 16587
$ws.Range("I103").Value = 16983.75
$ws.Range("J103").Value = 15000
$ws.Range("K103").Value = 16983.75
$ws.Range("L103").Value = 15000
$ws.Range("M103").Value = -15811.75
$ws.Range("N103").Value = -17344

$ws.Range("H107").Value = 1084.7778
$ws.Range("I107").Value = 483.4
$ws.Range("K107").Value = 483.4
$ws.Range("M107").Value = 1436.6

$ws.Range("H122").Value = 2995
$ws.Range("I122").Value = 2995
$ws.Range("K122").Value = 8985
$ws.Range("M122").Value = -6535

$ws.Range("H132").Value = 1625
$ws.Range("I132").Value = 1625
$ws.Range("K132").Value = 4875
$ws.Range("M132").Value = -2345

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 330
$ws.Range("I121").Value = 330
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 990
$ws.Range("L121").Value = 0
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = 320

$ws.Range("H131").Value = 1511.5
$ws.Range("I131").Value = 568.6923
$ws.Range("K131").Value = 1706.0769
$ws.Range("M131").Value = 3333.9231

$ws.Range("H137").Value = 5240
$ws.Range("J137").Value = 5367.5
$ws.Range("L137").Value = 16102.5
$ws.Range("N137").Value = -26302.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 7621.875
$ws.Range("I122").Value = 6994
$ws.Range("K122").Value = 20982
$ws.Range("M122").Value = -18532

$ws.Range("H132").Value = 1250
$ws.Range("I132").Value = 875
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 2625
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -95
$ws.Range("N132").Value = -11060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1599.2
$ws.Range("I16").Value = 665.3333
$ws.Range("J16").Value = 3000
$ws.Range("K16").Value = 665.3333
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = -495.3333
$ws.Range("N16").Value = -3340

$ws.Range("H40").Value = 5893.684
$ws.Range("I40").Value = 5792.706
$ws.Range("J40").Value = 6752
$ws.Range("K40").Value = 5792.706
$ws.Range("L40").Value = 6752
$ws.Range("M40").Value = -5656.706
$ws.Range("N40").Value = -7024

$ws.Range("H46").Value = 3829.8965
$ws.Range("I46").Value = 4635.5713
$ws.Range("J46").Value = 3573.5454
$ws.Range("K46").Value = 4635.5713
$ws.Range("L46").Value = 3573.5454
$ws.Range("M46").Value = -4447.5713
$ws.Range("N46").Value = -3949.5454

$ws.Range("H74").Value = 10000
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 10000
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws.Range("H96").Value = 70000
$ws.Range("J96").Value = 70000
$ws.Range("L96").Value = 70000
$ws.Range("N96").Value = -75492

$ws.Range("H100").Value = 8854.272000000001
$ws.Range("J100").Value = 9649.700000000001
$ws.Range("L100").Value = 9649.700000000001
$ws.Range("N100").Value = -10731.7

$ws.Range("H122").Value = 3749
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws.Range("H132").Value = 2050
$ws.Range("I132").Value = 2050
$ws.Range("K132").Value = 6150
$ws.Range("M132").Value = -3620

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 30000
$ws.Range("I75").Value = 30000
$ws.Range("K75").Value = 30000
$ws.Range("M75").Value = -29064

$ws.Range("H78").Value = 30000
$ws.Range("I78").Value = 30000
$ws.Range("K78").Value = 90000
$ws.Range("M78").Value = -85320

$ws.Range("H107").Value = 1000
$ws.Range("I107").Value = 1000
$ws.Range("K107").Value = 3000
$ws.Range("M107").Value = -1080

$ws.Range("H136").Value = 2362.1667
$ws.Range("I136").Value = 1423.4286
$ws.Range("K136").Value = 4270.2858
$ws.Range("M136").Value = -1720.2858
